$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D (date serial), J, K, L, M, P per row (2-15).
# Values are taken directly from the target OOXML (post-edit) state.
$rowsData = @{
    2  = @{ D = 44676; J = 120; K = 4000; L = 4500; M = 4250; P = 71 }
    3  = @{ D = 44627; J = 120; K = 4000; L = 4500; M = 4250; P = 71 }
    4  = @{ D = 44281; J = 120; K = 5500; L = 6000; M = 5750; P = 96 }
    5  = @{ D = 44242; J = 160; K = 5000; L = 5500; M = 5250; P = 88 }
    6  = @{ D = 44669; J = 130; K = 4500; L = 5000; M = 4750; P = 79 }
    7  = @{ D = 44648; J = 120; K = 6500; L = 7000; M = 6750; P = 112 }
    8  = @{ D = 44657; J = 100; K = 5000; L = 5500; M = 5250; P = 88 }
    9  = @{ D = 44589; J = 110; K = 5000; L = 6000; M = 5500; P = 92 }
    10 = @{ D = 44603; J = 140; K = 5500; L = 6000; M = 5750; P = 96 }
    11 = @{ D = 44400; J = 120; K = 9000; L = 10000; M = 9500; P = 158 }
    12 = @{ D = 44362; J = 120; K = 8000; L = 9000; M = 8500; P = 142 }
    13 = @{ D = 44382; J = 160; K = 7000; L = 8000; M = 7438; P = 124 }
    14 = @{ D = 44494; J = 120; K = 5000; L = 6000; M = 5500; P = 92 }
    15 = @{ D = 44421; J = 100; K = 8000; L = 9000; M = 8500; P = 142 }
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # column D
    $ws.Cells.Item($r, 10).Value = $vals.J   # column J
    $ws.Cells.Item($r, 11).Value = $vals.K   # column K
    $ws.Cells.Item($r, 12).Value = $vals.L   # column L
    $ws.Cells.Item($r, 13).Value = $vals.M   # column M
    $ws.Cells.Item($r, 16).Value = $vals.P   # column P
}
